$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proximity")

$rows = @(
    @{ Row = 26; A = "2026-02-01"; B = "18:17:37"; C = "18:00"; D = "Living Room Main Door"; E = "EXIT";  F = "User EXITED Living Room Main Door" },
    @{ Row = 27; A = "2026-02-01"; B = "18:18:41"; C = "18:00"; D = "Bathroom Door";          E = "ENTER"; F = "User ENTERED Bathroom" },
    @{ Row = 28; A = "2026-02-01"; B = "18:18:49"; C = "18:00"; D = "Bathroom Door";          E = "EXIT";  F = "User EXITED Bathroom" }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row

    # Column A holds a date-formatted string ("2026-02-01"). Excel's COM
    # layer auto-detects that pattern and stores it as a real date serial,
    # but the source log keeps it as literal text, so force text formatting
    # before assigning, then drop back to the Normal style so no stray
    # number-format is left attached to the cell.
    $cellA = $ws.Cells.Item($rowIndex, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.A
    $cellA.Style = "Normal"

    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
    $ws.Cells.Item($rowIndex, 5).Value = $r.E
    $ws.Cells.Item($rowIndex, 6).Value = $r.F
}
